$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$shp = $ws1.Shapes.AddChart2(-1, -4169)
$chart = $shp.Chart
$chart.SetSourceData($ws1.Range("A1:B5"))
$chart.ChartType = -4169
Write-Host "ok"
